$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2029.8928
$ws.Range("J17").Value = 2256.2273
$ws.Range("L17").Value = 6768.6819
$ws.Range("N17").Value = -7104.6819
$ws.Range("H69").Value = 8132.25
$ws.Range("I69").Value = 4999.6665
$ws.Range("K69").Value = 14998.9995
$ws.Range("M69").Value = -14124.9995
$ws.Range("H72").Value = 8132.25
$ws.Range("I72").Value = 4999.6665
$ws.Range("K72").Value = 44996.9985
$ws.Range("M72").Value = -40628.9985
$ws.Range("H137").Value = 1011357.9
$ws.Range("I137").Value = 556683.4399999999
$ws.Range("K137").Value = 1670050.32
$ws.Range("M137").Value = -1667500.32

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6950513.5
$ws.Range("I32").Value = 7048267
$ws.Range("J32").Value = 9999
$ws.Range("K32").Value = 7048267
$ws.Range("L32").Value = 9999
$ws.Range("M32").Value = -7047980
$ws.Range("N32").Value = -10573
$ws.Range("H45").Value = 7699.7144
$ws.Range("I45").Value = 7699.7144
$ws.Range("K45").Value = 7699.7144
$ws.Range("M45").Value = -7322.7144
$ws.Range("H61").Value = 527557.0600000001
$ws.Range("I61").Value = 958493.7
$ws.Range("J61").Value = 7461.1377
$ws.Range("K61").Value = 958493.7
$ws.Range("L61").Value = 7461.1377
$ws.Range("M61").Value = -958281.7
$ws.Range("N61").Value = -7885.1377
$ws.Range("H63").Value = 4134
$ws.Range("I63").Value = 2068.1428
$ws.Range("J63").Value = 7749.25
$ws.Range("K63").Value = 2068.1428
$ws.Range("L63").Value = 7749.25
$ws.Range("M63").Value = -1382.1428
$ws.Range("N63").Value = -9121.25
$ws.Range("H66").Value = 4134
$ws.Range("I66").Value = 2068.1428
$ws.Range("J66").Value = 7749.25
$ws.Range("K66").Value = 10340.714
$ws.Range("L66").Value = 38746.25
$ws.Range("M66").Value = -6908.714
$ws.Range("N66").Value = -45610.25
$ws.Range("H74").Value = 2275851.5
$ws.Range("I74").Value = 2979069
$ws.Range("J74").Value = 3918
$ws.Range("K74").Value = 2979069
$ws.Range("L74").Value = 3918
$ws.Range("M74").Value = -2978195
$ws.Range("N74").Value = -5666
$ws.Range("H77").Value = 2275851.5
$ws.Range("I77").Value = 2979069
$ws.Range("J77").Value = 3918
$ws.Range("K77").Value = 14895345
$ws.Range("L77").Value = 19590
$ws.Range("M77").Value = -14890977
$ws.Range("N77").Value = -28326
$ws.Range("H110").Value = 2178.125
$ws.Range("I110").Value = 2201.7144
$ws.Range("K110").Value = 2201.7144
$ws.Range("M110").Value = -156.7143999999998
$ws.Range("H132").Value = 249202.36
$ws.Range("I132").Value = 426123.25
$ws.Range("K132").Value = 1278369.75
$ws.Range("M132").Value = -1275839.75
$ws.Range("H136").Value = 527557.0600000001
$ws.Range("I136").Value = 958493.7
$ws.Range("J136").Value = 7461.1377
$ws.Range("K136").Value = 2875481.1
$ws.Range("L136").Value = 22383.4131
$ws.Range("M136").Value = -2872931.1
$ws.Range("N136").Value = -27483.4131

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 43671.3
$ws.Range("J88").Value = 43671.3
$ws.Range("L88").Value = 43671.3
$ws.Range("N88").Value = -44483.3
$ws.Range("H91").Value = 43671.3
$ws.Range("J91").Value = 43671.3
$ws.Range("L91").Value = 43671.3
$ws.Range("N91").Value = -46479.3
$ws.Range("H105").Value = 3480.5
$ws.Range("I105").Value = 2917.9285
$ws.Range("K105").Value = 2917.9285
$ws.Range("M105").Value = -1170.9285
$ws.Range("H107").Value = 1803.1111
$ws.Range("I107").Value = 2100.158
$ws.Range("K107").Value = 2100.158
$ws.Range("M107").Value = -180.1579999999999
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6934.827
$ws.Range("J31").Value = 9178.806
$ws.Range("L31").Value = 9178.806
$ws.Range("N31").Value = -9768.806
$ws.Range("H34").Value = 6934.827
$ws.Range("J34").Value = 9178.806
$ws.Range("L34").Value = 9178.806
$ws.Range("N34").Value = -9582.806
$ws.Range("H74").Value = 55636.3
$ws.Range("J74").Value = 53298.6
$ws.Range("L74").Value = 53298.6
$ws.Range("N74").Value = -55046.6
$ws.Range("H77").Value = 55636.3
$ws.Range("J77").Value = 53298.6
$ws.Range("L77").Value = 159895.8
$ws.Range("N77").Value = -168631.8
$ws.Range("H122").Value = 1933.6666
$ws.Range("I122").Value = 1960.4
$ws.Range("K122").Value = 5881.200000000001
$ws.Range("M122").Value = -3431.200000000001
$ws.Range("H132").Value = 5962419.5
$ws.Range("I132").Value = 10784.895
$ws.Range("J132").Value = 62502948
$ws.Range("K132").Value = 32354.685
$ws.Range("L132").Value = 187508844
$ws.Range("M132").Value = -29824.685
$ws.Range("N132").Value = -187513904
$ws.Range("H141").Value = 235900
$ws.Range("J141").Value = 248430.84
$ws.Range("L141").Value = 248430.84
$ws.Range("N141").Value = -258790.84

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11035.462
$ws.Range("J131").Value = 12037.728
$ws.Range("L131").Value = 36113.18399999999
$ws.Range("N131").Value = -46193.18399999999
$ws.Range("H132").Value = 2472.611
$ws.Range("I132").Value = 2172.0908
$ws.Range("K132").Value = 19548.8172
$ws.Range("M132").Value = -17018.8172
$ws.Range("H138").Value = 3897.8
$ws.Range("I138").Value = 3872.25
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 11616.75
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -6476.75
$ws.Range("N138").Value = -22280
$ws.Range("H139").Value = 906.0833
$ws.Range("I139").Value = 806.63635
$ws.Range("K139").Value = 2419.90905
$ws.Range("M139").Value = 2720.09095

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 728031.5
$ws.Range("I126").Value = 982161
$ws.Range("K126").Value = 2946483
$ws.Range("M126").Value = -2944013
$ws.Range("H134").Value = 51082.332
$ws.Range("J134").Value = 51635.363
$ws.Range("L134").Value = 154906.089
$ws.Range("N134").Value = -159976.089
$ws.Range("H136").Value = 37301.773
$ws.Range("J136").Value = 37301.773
$ws.Range("L136").Value = 111905.319
$ws.Range("N136").Value = -117005.319

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4844.375
$ws.Range("J7").Value = 3920
$ws.Range("L7").Value = 3920
$ws.Range("N7").Value = -4144
$ws.Range("H93").Value = 1844
$ws.Range("J93").Value = 1678.8
$ws.Range("L93").Value = 1678.8
$ws.Range("N93").Value = -4174.8
$ws.Range("H126").Value = 4844.375
$ws.Range("J126").Value = 3920
$ws.Range("L126").Value = 11760
$ws.Range("N126").Value = -16700
$ws.Range("H132").Value = 655261.5600000001
$ws.Range("I132").Value = 867074.5
$ws.Range("J132").Value = 3529.4614
$ws.Range("K132").Value = 2601223.5
$ws.Range("L132").Value = 10588.3842
$ws.Range("M132").Value = -2598693.5
$ws.Range("N132").Value = -15648.3842

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6806717.5
$ws.Range("I136").Value = 7329738.5
$ws.Range("K136").Value = 21989215.5
$ws.Range("M136").Value = -21986665.5
